# Model Comparison.xlsx edit:
#  - Insert a new "Income Dataset" worksheet (results of income classifiers)
#    between "Wine Dataset" and "Titanic Dataset".
#  - Update the saved selection / active cell on the Wine Dataset and
#    Titanic Dataset sheets.

$wb = $excel.ActiveWorkbook

$wine = $wb.Worksheets.Item("Wine Dataset")
$titanic = $wb.Worksheets.Item("Titanic Dataset")

# --- update the saved selection on the Wine Dataset sheet -----------------
$wine.Range("B13").Select() | Out-Null

# --- update the saved selection on the Titanic Dataset sheet --------------
$titanic.Range("C8").Select() | Out-Null

# --- create the new Income Dataset sheet, right before Titanic Dataset ----
$income = $wb.Worksheets.Add($titanic)
$income.Name = "Income Dataset"

# Text / label cells are entered column-by-column (A, then B, then C, then D)
# to match the order new shared-string entries were originally authored in.

# Column A - Model names
$income.Range("A1").Value = "Model"
$income.Range("A2").Value = "Decision Trees"
$income.Range("A3").Value = "Neural Networks"
$income.Range("A4").Value = "Ada Boost"
$income.Range("A5").Value = "Support Vector Machines"
$income.Range("A6").Value = "K-Nearest Neighbors"
$income.Range("A8").Value = "Guessing"

# Column B - Parameter1
$income.Range("B1").Value = "Parameter1"
$income.Range("B2").Value = "alpha: 0.06"
$income.Range("B3").Value = "alpha: 0.005"
$income.Range("B4").Value = "alpha: 0.06"
$income.Range("B5").Value = "kernel: linear"
$income.Range("B6").Value = "number of neighbors: 10"

# Column C - Parameter2
$income.Range("C1").Value = "Parameter2"
$income.Range("C3").Value = "hidden layer size: 3"
$income.Range("C4").Value = "learning rate: 0.01"
$income.Range("C5").Value = "gamma: 0.1"

# Column D - Parameter3
$income.Range("D1").Value = "Parameter3"
$income.Range("D4").Value = "number of estimators: 150"
$income.Range("D5").Value = "C: 10"

# Column E - Mean Score
$income.Range("E1").Value = "Mean Score"
$income.Range("E2").Value = 0.82499999999999996
$income.Range("E3").Value = 0.83
$income.Range("E4").Value = 0.84
$income.Range("E5").Value = 0.83799999999999997
$income.Range("E6").Value = 0.82299999999999995

# Column F - Training Set Accuracy
$income.Range("F1").Value = "Training Set Accuracy"
$income.Range("F2").Value = 0.83499999999999996
$income.Range("F3").Value = 0.86829999999999996
$income.Range("F4").Value = 0.86899999999999999
$income.Range("F5").Value = 0.85499999999999998
$income.Range("F6").Value = 0.82199999999999995
$income.Range("F8").Value = 0.5

# Column G - Test Set Accuracy
$income.Range("G1").Value = "Test Set Accuracy"
$income.Range("G2").Value = 0.83
$income.Range("G3").Value = 0.83540000000000003
$income.Range("G4").Value = 0.83299999999999996
$income.Range("G5").Value = 0.83899999999999997
$income.Range("G6").Value = 0.82699999999999996
$income.Range("G8").Value = 0.5

# Column H - Runtime in s
$income.Range("H1").Value = "Runtime in s"
$income.Range("H2").Value = 217
$income.Range("H3").Value = 2916
$income.Range("H4").Value = 7065
$income.Range("H5").Value = 271
$income.Range("H6").Value = 1456

# --- number formats --------------------------------------------------------
$percentFormat = "0.0%"
$runtimeFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

$income.Range("E1:G8").NumberFormat = $percentFormat
$income.Range("H1:H8").NumberFormat = $runtimeFormat

# --- column widths (characters) --------------------------------------------
$income.Columns.Item(1).ColumnWidth = 22.666666666666668
$income.Columns.Item(2).ColumnWidth = 21.166666666666668
$income.Columns.Item(3).ColumnWidth = 18.5
$income.Columns.Item(4).ColumnWidth = 23.666666666666668
$income.Columns.Item(5).ColumnWidth = 10.5
$income.Columns.Item(6).ColumnWidth = 19
$income.Columns.Item(7).ColumnWidth = 15.5
$income.Columns.Item(8).ColumnWidth = 12.666666666666666

# --- selection on the new sheet --------------------------------------------
$income.Range("B4").Select() | Out-Null
